# Revert "Prezentacija za Leap"
#
# 1) Remove the last two slides (the "Leap Motion kontroler" slide and the
#    "Upravljanje autom" slide) that were appended by the reverted commit.
# 2) Roll the cached "datetimeFigureOut" field text back from 1/8/2015 to
#    1/7/2015 everywhere it appears (slide master + every slide layout).

$p = $ppt.ActivePresentation

# --- 1. Drop the trailing two slides ----------------------------------
# They are the last two entries in the slide list (positions 10 and 11),
# delete from the end so indices of the earlier slides stay stable.
while ($p.Slides.Count -gt 9) {
    $p.Slides.Item($p.Slides.Count).Delete()
}

# --- 2. Fix up the cached date field text -----------------------------
function Update-DateFld($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "1/8/2015") {
                $tr.Text = "1/7/2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateFld $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateFld $layout.Shapes
}
